$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("E1").Value = "Execution Time (ms)"
$ws.Range("F1").Value = "Memory Usage (B)"

# Match the style of the existing header row (A1:D1) for the new headers
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-set values after paste (PasteSpecial formats only, but just to be safe)
$ws.Range("E1").Value = "Execution Time (ms)"
$ws.Range("F1").Value = "Memory Usage (B)"

# Execution Time (ms) values
$ws.Range("E2").Value = 7.519000006141141
$ws.Range("E3").Value = 4.476800007978454
$ws.Range("E4").Value = 42.64199998578988
$ws.Range("E5").Value = 3.086599987000227

# Memory Usage (B) values
$ws.Range("F2").Value = 32768
$ws.Range("F3").Value = 45056
$ws.Range("F4").Value = 8192
$ws.Range("F5").Value = 4096
